$d = $word.ActiveDocument

$d.Content.Find.Execute("27×57=", $true, $false, $false, $false, $false, $true, 1, $false, "71×59=", 2) | Out-Null
$d.Content.Find.Execute("94×32=", $true, $false, $false, $false, $false, $true, 1, $false, "92×57=", 2) | Out-Null
$d.Content.Find.Execute("19×32=", $true, $false, $false, $false, $false, $true, 1, $false, "98×83=", 2) | Out-Null
$d.Content.Find.Execute("26×30=", $true, $false, $false, $false, $false, $true, 1, $false, "44×70=", 2) | Out-Null
$d.Content.Find.Execute("37×43=", $true, $false, $false, $false, $false, $true, 1, $false, "74×72=", 2) | Out-Null
$d.Content.Find.Execute("88×30=", $true, $false, $false, $false, $false, $true, 1, $false, "73×89=", 2) | Out-Null
$d.Content.Find.Execute("73×70=", $true, $false, $false, $false, $false, $true, 1, $false, "70×16=", 2) | Out-Null
$d.Content.Find.Execute("43×66=", $true, $false, $false, $false, $false, $true, 1, $false, "41×30=", 2) | Out-Null
$d.Content.Find.Execute("41×54=", $true, $false, $false, $false, $false, $true, 1, $false, "27×29=", 2) | Out-Null
$d.Content.Find.Execute("78×35=", $true, $false, $false, $false, $false, $true, 1, $false, "86×26=", 2) | Out-Null
$d.Content.Find.Execute("81×82=", $true, $false, $false, $false, $false, $true, 1, $false, "82×72=", 2) | Out-Null
$d.Content.Find.Execute("12×86=", $true, $false, $false, $false, $false, $true, 1, $false, "64×65=", 2) | Out-Null
$d.Content.Find.Execute("91×71=", $true, $false, $false, $false, $false, $true, 1, $false, "85×67=", 2) | Out-Null
$d.Content.Find.Execute("74×93=", $true, $false, $false, $false, $false, $true, 1, $false, "88×89=", 2) | Out-Null
$d.Content.Find.Execute("48×51=", $true, $false, $false, $false, $false, $true, 1, $false, "49×47=", 2) | Out-Null
$d.Content.Find.Execute("92×93=", $true, $false, $false, $false, $false, $true, 1, $false, "94×61=", 2) | Out-Null
$d.Content.Find.Execute("18×54=", $true, $false, $false, $false, $false, $true, 1, $false, "86×40=", 2) | Out-Null
$d.Content.Find.Execute("86×41=", $true, $false, $false, $false, $false, $true, 1, $false, "84×41=", 2) | Out-Null
$d.Content.Find.Execute("16×13=", $true, $false, $false, $false, $false, $true, 1, $false, "74×86=", 2) | Out-Null
$d.Content.Find.Execute("20×30=", $true, $false, $false, $false, $false, $true, 1, $false, "60×72=", 2) | Out-Null
$d.Content.Find.Execute("91×49=", $true, $false, $false, $false, $false, $true, 1, $false, "40×39=", 2) | Out-Null
$d.Content.Find.Execute("34×26=", $true, $false, $false, $false, $false, $true, 1, $false, "52×17=", 2) | Out-Null
$d.Content.Find.Execute("99×77=", $true, $false, $false, $false, $false, $true, 1, $false, "77×90=", 2) | Out-Null
$d.Content.Find.Execute("74×21=", $true, $false, $false, $false, $false, $true, 1, $false, "33×26=", 2) | Out-Null
$d.Content.Find.Execute("70×18=", $true, $false, $false, $false, $false, $true, 1, $false, "96×91=", 2) | Out-Null
